$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..8) {
    $cell = $ws.Cells.Item($r, 3)  # Column C ("Förändrad")
    if ($cell.Value2 -eq 46060) {
        $cell.Value2 = 46061
    }
}
